# Publish a new version of the term (ValueSet-KLMatterOfInterestValues):
# bump Version 1.0.0 -> 1.1.0 and refresh the Date stamp, on the
# "Metadata" sheet (column A = Property name, column B = Value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$versionRow = $null
$dateRow = $null

for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2
    if ($label -eq "Version") { $versionRow = $r }
    if ($label -eq "Date") { $dateRow = $r }
}

if ($versionRow -ne $null) {
    $ws.Cells.Item($versionRow, 2).Value = "1.1.0"
}

if ($dateRow -ne $null) {
    $ws.Cells.Item($dateRow, 2).Value = "2023-07-10T23:08:03+02:00"
}
